$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a literal text value (avoids Excel's
# automatic number/date inference turning numeric-looking strings such
# as "1.003" or "10.08" into real numbers), while leaving the cell's
# style back at the workbook default ("Normal") once done so no stray
# custom number-format style is left behind.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

# Rows that only need Price (D) and Volume(1h) (E) updates
$rows = @{
    2  = @{ D = "30.251.31";     E = "  +1.00%  " }
    3  = @{ D = "1.919.32";      E = "  +0.57%  " }
    4  = @{ D = "1.003";         E = "  +0.02%  " }
    5  = @{ D = "0.8019";        E = "  +0.21%  " }
    6  = @{ D = "244.46";        E = "  +1.01%  " }
    7  = @{ D = "1.002";         E = "  +0.06%  " }
    8  = @{ D = "0.3255";        E = "  +3.37%  " }
    9  = @{ D = "26.88";         E = "  +2.37%  " }
    10 = @{ D = "0.07251";       E = "  +5.15%  " }
    11 = @{ D = "0.7863";        E = "  +6.91%  " }
    12 = @{ D = "0.08081";       E = "  +1.04%  " }
    13 = @{ D = "1.919.76";      E = "  +0.48%  " }
    14 = @{ D = "5.402";         E = "  +4.16%  " }
    15 = @{ D = "93.89";         E = "  +1.01%  " }
    16 = @{ D = "30.250.75";     E = "  +0.89%  " }
    17 = @{ D = "14.21";         E = "  +1.75%  " }
    18 = @{ D = "6.070";         E = "  +3.51%  " }
    19 = @{ D = "249.83";        E = "  +1.89%  " }
    20 = @{ D = "0.000007861";   E = "  +1.85%  " }
    23 = @{ E = "  +0.01%  " }
    24 = @{ D = "1.004";         E = "  +0.00%  " }
    25 = @{ D = "0.1626";        E = "  +14.63%  " }
    26 = @{ D = "9.466";         E = "  +2.79%  " }
    27 = @{ D = "167.61";        E = "  -0.02%  " }
    28 = @{ D = "19.00";         E = "  +0.46%  " }
    29 = @{ D = "2.152";         E = "  +6.13%  " }
    30 = @{ D = "1.392";         E = "  +2.25%  " }
    31 = @{ E = "  +2.51%  " }
    32 = @{ D = "4.497";         E = "  +4.54%  " }
    33 = @{ D = "0.05696";       E = "  +4.00%  " }
    34 = @{ D = "4.162";         E = "  +2.22%  " }
    35 = @{ D = "1.296";         E = "  +2.89%  " }
    36 = @{ D = "0.7510";        E = "  +2.48%  " }
    37 = @{ D = "1.003";         E = "  +0.30%  " }
    38 = @{ D = "2.732";         E = "  +0.36%  " }
    39 = @{ D = "0.01959";       E = "  +1.98%  " }
    40 = @{ E = "  +1.29%  " }
    41 = @{ D = "0.4525";        E = "  +2.57%  " }
    42 = @{ D = "74.05";         E = "  +2.58%  " }
    43 = @{ D = "6.018";         E = "  -2.51%  " }
    44 = @{ D = "0.8576";        E = "  +2.18%  " }
    45 = @{ D = "1.932";         E = "  +3.19%  " }
    46 = @{ E = "  +0.03%  " }
    47 = @{ D = "1.039.84";      E = "  +6.33%  " }
    48 = @{ D = "103.27";        E = "  +2.66%  " }
    51 = @{ D = "7.628";         E = "  +1.12%  " }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    if ($vals.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$r") $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$r").Value2 = $vals["E"]
    }
}

# Rows 21/22 swapped their coin (Chainlink <-> WrappedliquidstakedEther2.0) along with new data
$ws.Range("B21").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D21") "2.180.92"
$ws.Range("E21").Value2 = "  +0.26%  "

$ws.Range("B22").Value2 = "Chainlink"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D22") "8.226"
$ws.Range("E22").Value2 = "  +20.38%  "

# Rows 49/50 swapped their coin (SynthetixNetwork <-> EnergySwap) along with new data
$ws.Range("B49").Value2 = "EnergySwap"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "10.08"
$ws.Range("E49").Value2 = "  +3.59%  "

$ws.Range("B50").Value2 = "SynthetixNetwork"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D50") "3.108"
$ws.Range("E50").Value2 = "  +12.42%  "

Write-Output "done"
